# Update "想去人数" (interested-count) figures that were refreshed by the
# generator run at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 302
$wsExpo.Range("F4").Value = 1304

# Sheet "全部类型" (All types) - same two events appear here too
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 302
$wsAll.Range("F5").Value = 1304
